# Refresh script for cryptos.xlsx: updates Price (D) and Volume(1h) (E)
# columns with newly scraped values, and swaps the Maker / TrustWalletToken
# rows (41 <-> 42) to reflect their new rank ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values look like plain decimal numbers (e.g. "256.21").
# The source data stores every Price/Volume cell as text, so force those specific
# cells to a text number format before assigning them - this keeps them as strings
# instead of letting Excel auto-convert them to numeric values.
$textFormatCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D12",
    "D14",
    "D15",
    "D16",
    "D19",
    "D21",
    "D26",
    "D27",
    "D29",
    "D32",
    "D33",
    "D34",
    "D36",
    "D39",
    "D42",
    "D44",
    "D45",
    "D46",
    "D49",
    "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.695.47"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "1.991.40"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "256.21"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "55.35"
$ws.Range("E8").Value = "  -7.51%  "
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("E10").Value = "  -5.83%  "
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").Value = "14.20"
$ws.Range("E12").Value = "  -6.81%  "
$ws.Range("D13").Value = "2.285.23"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "21.31"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "0.790"
$ws.Range("E15").Value = "  -7.09%  "
$ws.Range("D16").Value = "5.16"
$ws.Range("E16").Value = "  -5.32%  "
$ws.Range("D17").Value = "1.979.46"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "36.664.56"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "70.46"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("D21").Value = "235.05"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "163.93"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "8.87"
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").Value = "1.35"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("E30").Value = "  -9.22%  "
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("D32").Value = "4.54"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "0.0628"
$ws.Range("E33").Value = "  -6.73%  "
$ws.Range("D34").Value = "4.35"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("E35").Value = "  -8.59%  "
$ws.Range("D36").Value = "3.47"
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "5.46"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.448.44"
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("D44").Value = "0.0909"
$ws.Range("E44").Value = "  -6.51%  "
$ws.Range("D45").Value = "88.50"
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "15.50"
$ws.Range("E46").Value = "  -7.17%  "
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "6.86"
$ws.Range("E49").Value = "  -9.37%  "
$ws.Range("D50").Value = "2.177.23"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  -9.39%  "
